$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Ephesians 1:1"
$ws.Cells.Item(3, 1).Value = "Ephesians 1:8"
$ws.Cells.Item(4, 1).Value = "Ephesians 1:10"
$ws.Cells.Item(5, 1).Value = "Ephesians 1:11"
$ws.Cells.Item(6, 1).Value = "Ephesians 1:12"
$ws.Cells.Item(7, 1).Value = "Ephesians 1:13"
$ws.Cells.Item(8, 1).Value = "Ephesians 1:14"
$ws.Cells.Item(9, 1).Value = "Ephesians 1:15"
$ws.Cells.Item(10, 1).Value = "Ephesians 1:16"
$ws.Cells.Item(11, 1).Value = "Ephesians 1:17"
$ws.Cells.Item(12, 1).Value = "Ephesians 1:19"
$ws.Cells.Item(13, 1).Value = "Ephesians 1:20"
$ws.Cells.Item(14, 1).Value = "Ephesians 1:22"
$ws.Cells.Item(15, 1).Value = "Ephesians 2:3"
$ws.Cells.Item(16, 1).Value = "Ephesians 2:4"
$ws.Cells.Item(17, 1).Value = "Ephesians 2:5"
$ws.Cells.Item(18, 1).Value = "Ephesians 2:6"
$ws.Cells.Item(19, 1).Value = "Ephesians 2:7"
$ws.Cells.Item(20, 1).Value = "Ephesians 2:8"
$ws.Cells.Item(21, 1).Value = "Ephesians 2:10"
$ws.Cells.Item(22, 1).Value = "Ephesians 2:11"
$ws.Cells.Item(24, 1).Value = "Ephesians 2:13"
$ws.Cells.Item(25, 1).Value = "Ephesians 2:14"
$ws.Cells.Item(26, 1).Value = "Ephesians 2:15"
$ws.Cells.Item(27, 1).Value = "Ephesians 2:16"
$ws.Cells.Item(28, 1).Value = "Ephesians 2:17"
$ws.Cells.Item(29, 1).Value = "Ephesians 2:18"
$ws.Cells.Item(30, 1).Value = "Ephesians 2:19"
$ws.Cells.Item(31, 1).Value = "Ephesians 2:22"
$ws.Cells.Item(32, 1).Value = "Ephesians 3:1"
$ws.Cells.Item(34, 1).Value = "Ephesians 3:5"
$ws.Cells.Item(35, 1).Value = "Ephesians 3:6"
$ws.Cells.Item(36, 1).Value = "Ephesians 3:8"
$ws.Cells.Item(37, 1).Value = "Ephesians 3:9"
$ws.Cells.Item(38, 1).Value = "Ephesians 3:10"
$ws.Cells.Item(39, 1).Value = "Ephesians 3:11"
$ws.Cells.Item(40, 1).Value = "Ephesians 3:12"
$ws.Cells.Item(41, 1).Value = "Ephesians 3:13"
$ws.Cells.Item(42, 1).Value = "Ephesians 3:14"
$ws.Cells.Item(43, 1).Value = "Ephesians 3:15"
$ws.Cells.Item(44, 1).Value = "Ephesians 3:16"
$ws.Cells.Item(46, 1).Value = "Ephesians 3:19"
$ws.Cells.Item(47, 1).Value = "Ephesians 4:1"
$ws.Cells.Item(48, 1).Value = "Ephesians 4:2"
$ws.Cells.Item(49, 1).Value = "Ephesians 4:3"
$ws.Cells.Item(50, 1).Value = "Ephesians 4:4"
$ws.Cells.Item(52, 1).Value = "Ephesians 4:7"
$ws.Cells.Item(53, 1).Value = "Ephesians 4:8"
$ws.Cells.Item(54, 1).Value = "Ephesians 4:9"
$ws.Cells.Item(55, 1).Value = "Ephesians 4:10"
$ws.Cells.Item(56, 1).Value = "Ephesians 4:11"
$ws.Cells.Item(57, 1).Value = "Ephesians 4:13"
$ws.Cells.Item(58, 1).Value = "Ephesians 4:14"
$ws.Cells.Item(59, 1).Value = "Ephesians 4:15"
$ws.Cells.Item(60, 1).Value = "Ephesians 4:17"
$ws.Cells.Item(61, 1).Value = "Ephesians 4:18"
$ws.Cells.Item(62, 1).Value = "Ephesians 4:19"
$ws.Cells.Item(63, 1).Value = "Ephesians 4:20"
$ws.Cells.Item(64, 1).Value = "Ephesians 4:21"
$ws.Cells.Item(65, 1).Value = "Ephesians 4:22"
$ws.Cells.Item(66, 1).Value = "Ephesians 4:23"
$ws.Cells.Item(67, 1).Value = "Ephesians 4:24"
$ws.Cells.Item(68, 1).Value = "Ephesians 4:25"
$ws.Cells.Item(69, 1).Value = "Ephesians 4:27"
$ws.Cells.Item(70, 1).Value = "Ephesians 4:28"
$ws.Cells.Item(71, 1).Value = "Ephesians 4:30"
$ws.Cells.Item(72, 1).Value = "Ephesians 5:2"
$ws.Cells.Item(73, 1).Value = "Ephesians 5:3"
$ws.Cells.Item(74, 1).Value = "Ephesians 5:4"
$ws.Cells.Item(75, 1).Value = "Ephesians 5:5"
$ws.Cells.Item(76, 1).Value = "Ephesians 5:7"
$ws.Cells.Item(77, 1).Value = "Ephesians 5:8"
$ws.Cells.Item(78, 1).Value = "Ephesians 5:9"
$ws.Cells.Item(79, 1).Value = "Ephesians 5:10"
$ws.Cells.Item(80, 1).Value = "Ephesians 5:12"
$ws.Cells.Item(81, 1).Value = "Ephesians 5:13"
$ws.Cells.Item(82, 1).Value = "Ephesians 5:14"
$ws.Cells.Item(83, 1).Value = "Ephesians 5:17"
$ws.Cells.Item(84, 1).Value = "Ephesians 5:18"
$ws.Cells.Item(85, 1).Value = "Ephesians 5:19"
$ws.Cells.Item(86, 1).Value = "Ephesians 5:23"
$ws.Cells.Item(87, 1).Value = "Ephesians 5:24"
$ws.Cells.Item(88, 1).Value = "Ephesians 5:26"
$ws.Cells.Item(89, 1).Value = "Ephesians 5:27"
$ws.Cells.Item(90, 1).Value = "Ephesians 5:28"
$ws.Cells.Item(91, 1).Value = "Ephesians 5:29"
$ws.Cells.Item(92, 1).Value = "Ephesians 5:30"
$ws.Cells.Item(93, 1).Value = "Ephesians 5:32"
$ws.Cells.Item(94, 1).Value = "Ephesians 5:33"
$ws.Cells.Item(95, 1).Value = "Ephesians 6:3"
$ws.Cells.Item(96, 1).Value = "Ephesians 6:4"
$ws.Cells.Item(97, 1).Value = "Ephesians 6:6"
$ws.Cells.Item(98, 1).Value = "Ephesians 6:7"
$ws.Cells.Item(99, 1).Value = "Ephesians 6:9"
$ws.Cells.Item(100, 1).Value = "Ephesians 6:10"
$ws.Cells.Item(101, 1).Value = "Ephesians 6:12"
$ws.Cells.Item(102, 1).Value = "Ephesians 6:13"
$ws.Cells.Item(103, 1).Value = "Ephesians 6:14"
$ws.Cells.Item(104, 1).Value = "Ephesians 6:15"
$ws.Cells.Item(105, 1).Value = "Ephesians 6:16"
$ws.Cells.Item(106, 1).Value = "Ephesians 6:17"
$ws.Cells.Item(107, 1).Value = "Ephesians 6:18"
$ws.Cells.Item(108, 1).Value = "Ephesians 6:19"
$ws.Cells.Item(109, 1).Value = "Ephesians 6:20"
$ws.Cells.Item(110, 1).Value = "Ephesians 6:22"
$ws.Cells.Item(111, 1).Value = "Ephesians 6:24"
